# Add season record columns (Wins / Losses / Ties) to the team stats sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from the last existing header cell (AC1) onto the
# three new header cells so they pick up the same bold/border/alignment
# formatting used by the rest of row 1.
$ws.Range("AC1").Copy($ws.Range("AD1:AF1"))

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$wins = 73
$losses = 89
$ties = 0

$lastRow = 43
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins    # column AD
    $ws.Cells.Item($r, 31).Value = $losses  # column AE
    $ws.Cells.Item($r, 32).Value = $ties    # column AF
}
